$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quick Notes ")

# Select the sheet so the saved view reflects it as the active sheet (tabSelected).
$ws.Select()

# --- Preserve hyperlinks across the upcoming row delete -------------------
# Deleting row 1 does not renumber the worksheet's Hyperlinks collection in
# this runtime, and re-adding a hyperlink via Hyperlinks.Add() overwrites the
# cell's existing direct formatting with the default "Hyperlink" look. To
# avoid both problems: stash each hyperlinked cell's exact formatting (via a
# format-only copy/paste, which reuses the existing style record instead of
# fabricating a new one) into a scratch cell far outside the used range,
# then restore it onto the relocated cell after the hyperlink is recreated.
$scratchCol = 50
$links = @()
$idx = 0
foreach ($hl in $ws.Hyperlinks) {
    $srcCell = $hl.Range
    $scratch = $ws.Cells.Item($srcCell.Row, $scratchCol + $idx)
    $srcCell.Copy()
    $scratch.PasteSpecial(-4122) # xlPasteFormats
    $links += , @($srcCell.Row, $srcCell.Column, $hl.Address, $scratch.Row, $scratch.Column)
    $idx = $idx + 1
}

# Delete the first row (the blank title row), shifting everything else up by one.
$ws.Rows.Item(1).Delete()

# Rebuild the hyperlinks at their correctly-shifted (row - 1) locations, then
# restore the original formatting that Hyperlinks.Add() just clobbered.
$ws.Hyperlinks.Delete()
foreach ($l in $links) {
    $newRow = $l[0] - 1
    $cell = $ws.Cells.Item($newRow, $l[1])
    $ws.Hyperlinks.Add($cell, $l[2])

    # The scratch cells live on rows outside the touched columns, but the
    # row-1 delete shifted every row (including the scratch rows) up by one.
    $scratchRow = $l[3] - 1
    $scratch = $ws.Cells.Item($scratchRow, $l[4])
    $scratch.Copy()
    $cell.PasteSpecial(-4122) # xlPasteFormats
}

# Clean up the scratch area so it doesn't widen the sheet's used range.
$ws.Application.CutCopyMode = $false
$ws.Range($ws.Cells.Item(1, $scratchCol), $ws.Cells.Item(20, $scratchCol + 10)).Clear()

# Reset the view so the top-left cell is A1 and the whole (new) first row is
# selected, matching what Excel does after deleting row 1 while it was
# visible at the top of the window.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Rows.Item(1).EntireRow.Select()
